$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 1.75
$ws.Range("I3").Value = 4
$ws.Range("AJ3").Value = 26

# Row 5
$ws.Range("G5").Value = 2.5
$ws.Range("I5").Value = 2.8
$ws.Range("W5").Value = 23
$ws.Range("AG5").Value = 11
$ws.Range("AI5").Value = 23
$ws.Range("AJ5").Value = 34

# Row 17
$ws.Range("G17").Value = 3.1
$ws.Range("H17").Value = 3
$ws.Range("T17").Value = 8.75
$ws.Range("U17").Value = 16
$ws.Range("W17").Value = 40
$ws.Range("X17").Value = 28
$ws.Range("AA17").Value = 5.8
$ws.Range("AF17").Value = 11.25
$ws.Range("AI17").Value = 21
$ws.Range("AJ17").Value = 32

# Row 24
$ws.Range("K24").Value = 10
